# QandAData.xlsx edit: add a helper column "y" (row index 0..28) as column C,
# bump a couple of indent levels, resize columns A/B, adjust row 28 height,
# and move the active selection to B18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column C: header "y" then sequential numbers 0..28 ---
$ws.Range("C1").Value = "y"
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 3).Value = $r - 2
}

# --- Indent bumps on the existing B-column styles ---
# B6  (indent 4  -> 7), keep its wrap-text on
$ws.Range("B6").WrapText = $true
$ws.Range("B6").IndentLevel = 7

# B8  (indent 6  -> 11)
$ws.Range("B8").IndentLevel = 11

# B13 (indent 9  -> 15)
$ws.Range("B13").IndentLevel = 15

# B17 and B19 (indent 9 -> 15)
$ws.Range("B17").IndentLevel = 15
$ws.Range("B19").IndentLevel = 15

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 44.9919028340081
$ws.Columns.Item(2).ColumnWidth = 84.1943319838057

# --- Row 28 height 45 -> 60 ---
$ws.Rows.Item(28).RowHeight = 60

# --- Move the selection / view to B18 ---
$ws.Range("B18").Select()
